$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Added the build forge and cannon strat so the ai can defend early game"
#    -> "...cannon start so the ai..."  (fix misspelling; drop the now
#    unnecessary spell-check proofErr wrapper around the word)
# ---------------------------------------------------------------------------
$para16Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Added the build forge and cannon </w:t></w:r><w:r><w:t>start</w:t></w:r><w:r><w:t xml:space="preserve"> so the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> can defend early game</w:t></w:r></w:p>
"@

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Added the build forge and cannon strat so the ai*") {
        $p.Range.InsertXML($para16Xml) | Out-Null
        $found = $true
        break
    }
}
if (-not $found) {
    throw "Could not locate the 'Added the build forge...' paragraph"
}

# ---------------------------------------------------------------------------
# 2) "IVE GIVE" + bookmark + "N UP" -> "IVE GIVEN UP", then append the rest
#    of the diary: the 20/02/2018 and 22/02/2018 entries, finishing with an
#    empty paragraph that now just carries the _GoBack bookmark.
# ---------------------------------------------------------------------------
$restXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>IVE GIVEN UP</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>20/02/2018</w:t></w:r></w:p><w:p><w:r><w:t>After meeting with my supervisor for most of the day, he helped me get back on track with the bot. All I need to do now is figure out a way to get the probes to build a little further away from the base and then that should solve the issues. Though a new issue cropped up with the scouting breaking the bot, but with a little reverting of code this fixed the issue, it was strange as I had done nothing to the behaviours of the scouting, I had only changed the plan’s code, yet this was enough to break it. At the end of this day I did not manage to get the probe to build elsewhere.</w:t></w:r></w:p><w:p><w:r><w:t>22/02/2018</w:t></w:r></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Today I have managed to get the chokepoint tile position, this has allowed me to give the builder a location to build, unfortunately the build command in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bwapi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> keeps returning false, and I’m not sure why as I’m passing in a tile position, perhaps its something to do with the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>PossibleBuildLocation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> function.</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@

$found2 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "IVE GIVEN UP" -or $t -eq "IVE GIVE") {
        $p.Range.InsertXML($restXml) | Out-Null
        $found2 = $true
        break
    }
}
if (-not $found2) {
    throw "Could not locate the 'IVE GIVE...' paragraph"
}
